$wb = $excel.ActiveWorkbook

# --- MCF sheet: hydro row max capacity factor changes from 0 to 0.95 ---
$mcf = $wb.Worksheets.Item("MCF")
$mcf.Range("B6").Value = 0.95

# --- About sheet: clarify that hydro is included in the non-variable plant note ---
$about = $wb.Worksheets.Item("About")
$about.Range("A10").Value = "in any given hour. This is used for non-variable plant types, including hydro. We apply"

# --- cursor/selection positions (cosmetic, matches author's final view state) ---
$about.Activate()
$about.Range("A11").Select()

$mcf.Activate()
$mcf.Range("B7").Select()

$about.Activate()
